# Importação dos dados na planilha excel
# Adiciona uma nova linha (linha 7) com os dados importados:
#   A7 e B7 ficam vazios, C7 recebe "UNI" e D7 recebe uma quebra de linha.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 1).Value = ""
$ws.Cells.Item(7, 2).Value = ""
$ws.Cells.Item(7, 3).Value = "UNI"
$ws.Cells.Item(7, 4).Value = "`n"
